$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, [string]$text) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '30.134.97'
Set-TextValue $ws.Range('E2') '  +3.00%  '
Set-TextValue $ws.Range('D3') '1.894.51'
Set-TextValue $ws.Range('E3') '  +0.00%  '
Set-TextValue $ws.Range('D5') '325.80'
Set-TextValue $ws.Range('E5') '  +3.37%  '
Set-TextValue $ws.Range('E6') '  -0.21%  '
Set-TextValue $ws.Range('D7') '0.5173'
Set-TextValue $ws.Range('E7') '  +0.54%  '
Set-TextValue $ws.Range('D8') '0.4007'
Set-TextValue $ws.Range('E8') '  +2.27%  '
Set-TextValue $ws.Range('D9') '0.08432'
Set-TextValue $ws.Range('E9') '  +0.19%  '
Set-TextValue $ws.Range('D10') '42.74'
Set-TextValue $ws.Range('E10') '  +0.67%  '
Set-TextValue $ws.Range('D11') '1.118'
Set-TextValue $ws.Range('E11') '  +0.44%  '
Set-TextValue $ws.Range('D12') '23.17'
Set-TextValue $ws.Range('E12') '  +12.32%  '
Set-TextValue $ws.Range('D13') '6.436'
Set-TextValue $ws.Range('E13') '  +3.19%  '
Set-TextValue $ws.Range('D14') '1.897.97'
Set-TextValue $ws.Range('E14') '  +0.12%  '
Set-TextValue $ws.Range('D15') '7.332'
Set-TextValue $ws.Range('E15') '  +0.36%  '
Set-TextValue $ws.Range('D16') '1.004'
Set-TextValue $ws.Range('E16') '  -0.16%  '
Set-TextValue $ws.Range('D17') '94.65'
Set-TextValue $ws.Range('E17') '  +1.74%  '
Set-TextValue $ws.Range('D18') '0.00001111'
Set-TextValue $ws.Range('E18') '  +0.53%  '
Set-TextValue $ws.Range('D19') '0.06664'
Set-TextValue $ws.Range('E19') '  -1.15%  '
Set-TextValue $ws.Range('D20') '18.27'
Set-TextValue $ws.Range('D21') '1.002'
Set-TextValue $ws.Range('E21') '  -0.19%  '
Set-TextValue $ws.Range('E22') '  -1.02%  '
Set-TextValue $ws.Range('D23') '30.153.43'
Set-TextValue $ws.Range('E23') '  +3.00%  '
Set-TextValue $ws.Range('D24') '11.29'
Set-TextValue $ws.Range('E24') '  +1.51%  '
Set-TextValue $ws.Range('D25') '2.216'
Set-TextValue $ws.Range('E25') '  +0.03%  '
Set-TextValue $ws.Range('D26') '21.94'
Set-TextValue $ws.Range('E26') '  +4.79%  '
Set-TextValue $ws.Range('D27') '2.117.32'
Set-TextValue $ws.Range('E27') '  +0.23%  '
Set-TextValue $ws.Range('D28') '161.25'
Set-TextValue $ws.Range('E28') '  +1.41%  '
Set-TextValue $ws.Range('D29') '2.384'
Set-TextValue $ws.Range('E29') '  -2.00%  '
Set-TextValue $ws.Range('D30') '128.94'
Set-TextValue $ws.Range('E30') '  +0.86%  '
Set-TextValue $ws.Range('E31') '  +3.45%  '
Set-TextValue $ws.Range('D32') '0.1056'
Set-TextValue $ws.Range('E32') '  +0.98%  '
Set-TextValue $ws.Range('D33') '6.089'
Set-TextValue $ws.Range('E33') '  -0.42%  '
Set-TextValue $ws.Range('D34') '3.700'
Set-TextValue $ws.Range('E34') '  +1.26%  '
Set-TextValue $ws.Range('D35') '0.02493'
Set-TextValue $ws.Range('E35') '  +0.77%  '
Set-TextValue $ws.Range('D36') '0.06554'
Set-TextValue $ws.Range('E36') '  +0.29%  '
Set-TextValue $ws.Range('D37') '0.2208'
Set-TextValue $ws.Range('E37') '  +0.87%  '
Set-TextValue $ws.Range('D38') '5.246'
Set-TextValue $ws.Range('E38') '  +2.42%  '
Set-TextValue $ws.Range('E39') '  -0.76%  '
Set-TextValue $ws.Range('D40') '11.78'
Set-TextValue $ws.Range('E40') '  +4.87%  '
Set-TextValue $ws.Range('D41') '8.764'
Set-TextValue $ws.Range('E41') '  -3.07%  '
Set-TextValue $ws.Range('D42') '0.6499'
Set-TextValue $ws.Range('E42') '  +0.10%  '
Set-TextValue $ws.Range('D43') '1.236'
Set-TextValue $ws.Range('E43') '  +0.35%  '
Set-TextValue $ws.Range('D44') '0.6104'
Set-TextValue $ws.Range('E44') '  +0.90%  '
Set-TextValue $ws.Range('D45') '13.28'
Set-TextValue $ws.Range('E45') '  +0.59%  '
Set-TextValue $ws.Range('E46') '  +0.84%  '
Set-TextValue $ws.Range('D47') '2.054'
Set-TextValue $ws.Range('E47') '  +0.67%  '
Set-TextValue $ws.Range('D48') '1.238'
Set-TextValue $ws.Range('E48') '  +0.79%  '
Set-TextValue $ws.Range('D49') '124.41'
Set-TextValue $ws.Range('E49') '  +1.17%  '
Set-TextValue $ws.Range('E50') '  -1.31%  '
Set-TextValue $ws.Range('D51') '79.10'
Set-TextValue $ws.Range('E51') '  +2.04%  '
